$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng_2_0 = $ws.Range("D2:G2")
$rng_2_0.NumberFormat = "@"
$ws.Range("D2").Value = "306.53"
$ws.Range("E2").Value = "-6.65%"
$ws.Range("F2").Value = "10-2-2023"
$ws.Range("G2").Value = "0"
$rng_2_0.ClearFormats()

$rng_3_0 = $ws.Range("D3:G3")
$rng_3_0.NumberFormat = "@"
$ws.Range("D3").Value = "39.29"
$ws.Range("E3").Value = "-13.00%"
$ws.Range("F3").Value = "10-2-2023"
$ws.Range("G3").Value = "0"
$rng_3_0.ClearFormats()

$rng_4_0 = $ws.Range("D4:G4")
$rng_4_0.NumberFormat = "@"
$ws.Range("D4").Value = "4.987"
$ws.Range("E4").Value = "-5.74%"
$ws.Range("F4").Value = "10-2-2023"
$ws.Range("G4").Value = "0"
$rng_4_0.ClearFormats()

$rng_5_0 = $ws.Range("D5:G5")
$rng_5_0.NumberFormat = "@"
$ws.Range("D5").Value = "0.07716"
$ws.Range("E5").Value = "-7.90%"
$ws.Range("F5").Value = "10-2-2023"
$ws.Range("G5").Value = "0"
$rng_5_0.ClearFormats()

$rng_6_0 = $ws.Range("D6:G6")
$rng_6_0.NumberFormat = "@"
$ws.Range("D6").Value = "4.264"
$ws.Range("E6").Value = "-3.69%"
$ws.Range("F6").Value = "10-2-2023"
$ws.Range("G6").Value = "0"
$rng_6_0.ClearFormats()

$rng_7_0 = $ws.Range("D7:G7")
$rng_7_0.NumberFormat = "@"
$ws.Range("D7").Value = "1.589"
$ws.Range("E7").Value = "-18.77%"
$ws.Range("F7").Value = "10-2-2023"
$ws.Range("G7").Value = "0"
$rng_7_0.ClearFormats()

$rng_8_0 = $ws.Range("D8:G8")
$rng_8_0.NumberFormat = "@"
$ws.Range("D8").Value = "0.9087"
$ws.Range("E8").Value = "-6.62%"
$ws.Range("F8").Value = "10-2-2023"
$ws.Range("G8").Value = "0"
$rng_8_0.ClearFormats()

$rng_9_0 = $ws.Range("E9:G9")
$rng_9_0.NumberFormat = "@"
$ws.Range("E9").Value = "-8.60%"
$ws.Range("F9").Value = "10-2-2023"
$ws.Range("G9").Value = "0"
$rng_9_0.ClearFormats()

$rng_10_0 = $ws.Range("D10:G10")
$rng_10_0.NumberFormat = "@"
$ws.Range("D10").Value = "0.1723"
$ws.Range("E10").Value = "-10.14%"
$ws.Range("F10").Value = "10-2-2023"
$ws.Range("G10").Value = "0"
$rng_10_0.ClearFormats()

$rng_11_0 = $ws.Range("D11:G11")
$rng_11_0.NumberFormat = "@"
$ws.Range("D11").Value = "0.08965"
$ws.Range("E11").Value = "-7.57%"
$ws.Range("F11").Value = "10-2-2023"
$ws.Range("G11").Value = "0"
$rng_11_0.ClearFormats()

$rng_12_0 = $ws.Range("D12:G12")
$rng_12_0.NumberFormat = "@"
$ws.Range("D12").Value = "0.04437"
$ws.Range("E12").Value = "-3.84%"
$ws.Range("F12").Value = "10-2-2023"
$ws.Range("G12").Value = "0"
$rng_12_0.ClearFormats()

$rng_13_0 = $ws.Range("D13:G13")
$rng_13_0.NumberFormat = "@"
$ws.Range("D13").Value = "7.041"
$ws.Range("E13").Value = "-15.61%"
$ws.Range("F13").Value = "10-2-2023"
$ws.Range("G13").Value = "0"
$rng_13_0.ClearFormats()

$rng_14_0 = $ws.Range("D14:G14")
$rng_14_0.NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").Value = "-0.18%"
$ws.Range("F14").Value = "10-2-2023"
$ws.Range("G14").Value = "0"
$rng_14_0.ClearFormats()

$rng_15_0 = $ws.Range("D15")
$rng_15_0.NumberFormat = "@"
$rng_15_1 = $ws.Range("F15:G15")
$rng_15_1.NumberFormat = "@"
$ws.Range("D15").Value = "0.001253"
$ws.Range("F15").Value = "10-2-2023"
$ws.Range("G15").Value = "0"
$rng_15_0.ClearFormats()
$rng_15_1.ClearFormats()

$rng_16_0 = $ws.Range("D16:G16")
$rng_16_0.NumberFormat = "@"
$ws.Range("D16").Value = "0.005652"
$ws.Range("E16").Value = "-1.63%"
$ws.Range("F16").Value = "10-2-2023"
$ws.Range("G16").Value = "0"
$rng_16_0.ClearFormats()

$rng_17_0 = $ws.Range("B17:G17")
$rng_17_0.NumberFormat = "@"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.361"
$ws.Range("E17").Value = "-0.06%"
$ws.Range("F17").Value = "10-2-2023"
$ws.Range("G17").Value = "0"
$rng_17_0.ClearFormats()

$rng_18_0 = $ws.Range("B18:G18")
$rng_18_0.NumberFormat = "@"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "2.591"
$ws.Range("E18").Value = "2.65%"
$ws.Range("F18").Value = "10-2-2023"
$ws.Range("G18").Value = "0"
$rng_18_0.ClearFormats()

$rng_19_0 = $ws.Range("B19:G19")
$rng_19_0.NumberFormat = "@"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "0.3366"
$ws.Range("E19").Value = "0.32%"
$ws.Range("F19").Value = "10-2-2023"
$ws.Range("G19").Value = "0"
$rng_19_0.ClearFormats()

$rng_20_0 = $ws.Range("B20:G20")
$rng_20_0.NumberFormat = "@"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "0.1343"
$ws.Range("E20").Value = "-0.62%"
$ws.Range("F20").Value = "10-2-2023"
$ws.Range("G20").Value = "0"
$rng_20_0.ClearFormats()

$rng_21_0 = $ws.Range("B21:G21")
$rng_21_0.NumberFormat = "@"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "0.2860"
$ws.Range("E21").Value = "5.16%"
$ws.Range("F21").Value = "10-2-2023"
$ws.Range("G21").Value = "0"
$rng_21_0.ClearFormats()

$rng_22_0 = $ws.Range("B22:G22")
$rng_22_0.NumberFormat = "@"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "0.04140"
$ws.Range("E22").Value = "-0.95%"
$ws.Range("F22").Value = "10-2-2023"
$ws.Range("G22").Value = "0"
$rng_22_0.ClearFormats()

$rng_23_0 = $ws.Range("B23:G23")
$rng_23_0.NumberFormat = "@"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D23").Value = "0.001192"
$ws.Range("E23").Value = "-3.60%"
$ws.Range("F23").Value = "10-2-2023"
$ws.Range("G23").Value = "0"
$rng_23_0.ClearFormats()

$rng_24_0 = $ws.Range("B24:G24")
$rng_24_0.NumberFormat = "@"
$ws.Range("B24").Value = "HotbitToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D24").Value = "0.004075"
$ws.Range("E24").Value = "-9.32%"
$ws.Range("F24").Value = "10-2-2023"
$ws.Range("G24").Value = "0"
$rng_24_0.ClearFormats()

$rng_25_0 = $ws.Range("D25:G25")
$rng_25_0.NumberFormat = "@"
$ws.Range("D25").Value = "0.0001225"
$ws.Range("E25").Value = "-5.74%"
$ws.Range("F25").Value = "10-2-2023"
$ws.Range("G25").Value = "0"
$rng_25_0.ClearFormats()

$rng_26_0 = $ws.Range("D26:G26")
$rng_26_0.NumberFormat = "@"
$ws.Range("D26").Value = "0.0002992"
$ws.Range("E26").Value = "0.38%"
$ws.Range("F26").Value = "10-2-2023"
$ws.Range("G26").Value = "0"
$rng_26_0.ClearFormats()

$rng_27_0 = $ws.Range("F27:G27")
$rng_27_0.NumberFormat = "@"
$ws.Range("F27").Value = "10-2-2023"
$ws.Range("G27").Value = "0"
$rng_27_0.ClearFormats()

$rng_28_0 = $ws.Range("F28:G28")
$rng_28_0.NumberFormat = "@"
$ws.Range("F28").Value = "10-2-2023"
$ws.Range("G28").Value = "0"
$rng_28_0.ClearFormats()

$rng_29_0 = $ws.Range("F29:G29")
$rng_29_0.NumberFormat = "@"
$ws.Range("F29").Value = "10-2-2023"
$ws.Range("G29").Value = "0"
$rng_29_0.ClearFormats()

$rng_30_0 = $ws.Range("F30:G30")
$rng_30_0.NumberFormat = "@"
$ws.Range("F30").Value = "10-2-2023"
$ws.Range("G30").Value = "0"
$rng_30_0.ClearFormats()

$rng_31_0 = $ws.Range("F31:G31")
$rng_31_0.NumberFormat = "@"
$ws.Range("F31").Value = "10-2-2023"
$ws.Range("G31").Value = "0"
$rng_31_0.ClearFormats()

$rng_32_0 = $ws.Range("F32:G32")
$rng_32_0.NumberFormat = "@"
$ws.Range("F32").Value = "10-2-2023"
$ws.Range("G32").Value = "0"
$rng_32_0.ClearFormats()

$rng_33_0 = $ws.Range("F33:G33")
$rng_33_0.NumberFormat = "@"
$ws.Range("F33").Value = "10-2-2023"
$ws.Range("G33").Value = "0"
$rng_33_0.ClearFormats()

$rng_34_0 = $ws.Range("F34:G34")
$rng_34_0.NumberFormat = "@"
$ws.Range("F34").Value = "10-2-2023"
$ws.Range("G34").Value = "0"
$rng_34_0.ClearFormats()

$rng_35_0 = $ws.Range("F35:G35")
$rng_35_0.NumberFormat = "@"
$ws.Range("F35").Value = "10-2-2023"
$ws.Range("G35").Value = "0"
$rng_35_0.ClearFormats()

$rng_36_0 = $ws.Range("F36:G36")
$rng_36_0.NumberFormat = "@"
$ws.Range("F36").Value = "10-2-2023"
$ws.Range("G36").Value = "0"
$rng_36_0.ClearFormats()

$rng_37_0 = $ws.Range("F37:G37")
$rng_37_0.NumberFormat = "@"
$ws.Range("F37").Value = "10-2-2023"
$ws.Range("G37").Value = "0"
$rng_37_0.ClearFormats()

$rng_38_0 = $ws.Range("D38:G38")
$rng_38_0.NumberFormat = "@"
$ws.Range("D38").Value = "0.02342"
$ws.Range("E38").Value = "-13.48%"
$ws.Range("F38").Value = "10-2-2023"
$ws.Range("G38").Value = "0"
$rng_38_0.ClearFormats()

$rng_39_0 = $ws.Range("D39:G39")
$rng_39_0.NumberFormat = "@"
$ws.Range("D39").Value = "0.05116"
$ws.Range("E39").Value = "-9.35%"
$ws.Range("F39").Value = "10-2-2023"
$ws.Range("G39").Value = "0"
$rng_39_0.ClearFormats()

$rng_40_0 = $ws.Range("E40:G40")
$rng_40_0.NumberFormat = "@"
$ws.Range("E40").Value = "2.18%"
$ws.Range("F40").Value = "10-2-2023"
$ws.Range("G40").Value = "0"
$rng_40_0.ClearFormats()

$rng_41_0 = $ws.Range("E41:G41")
$rng_41_0.NumberFormat = "@"
$ws.Range("E41").Value = "-6.11%"
$ws.Range("F41").Value = "10-2-2023"
$ws.Range("G41").Value = "0"
$rng_41_0.ClearFormats()

$rng_42_0 = $ws.Range("D42:G42")
$rng_42_0.NumberFormat = "@"
$ws.Range("D42").Value = "0.007614"
$ws.Range("E42").Value = "3.69%"
$ws.Range("F42").Value = "10-2-2023"
$ws.Range("G42").Value = "0"
$rng_42_0.ClearFormats()

$rng_43_0 = $ws.Range("D43:G43")
$rng_43_0.NumberFormat = "@"
$ws.Range("D43").Value = "0.002038"
$ws.Range("E43").Value = "-0.09%"
$ws.Range("F43").Value = "10-2-2023"
$ws.Range("G43").Value = "0"
$rng_43_0.ClearFormats()

$rng_44_0 = $ws.Range("D44:G44")
$rng_44_0.NumberFormat = "@"
$ws.Range("D44").Value = "0.008022"
$ws.Range("E44").Value = "1.29%"
$ws.Range("F44").Value = "10-2-2023"
$ws.Range("G44").Value = "0"
$rng_44_0.ClearFormats()

$rng_45_0 = $ws.Range("D45:G45")
$rng_45_0.NumberFormat = "@"
$ws.Range("D45").Value = "0.3310"
$ws.Range("E45").Value = "-5.83%"
$ws.Range("F45").Value = "10-2-2023"
$ws.Range("G45").Value = "0"
$rng_45_0.ClearFormats()

$rng_46_0 = $ws.Range("D46:G46")
$rng_46_0.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006685"
$ws.Range("E46").Value = "-4.29%"
$ws.Range("F46").Value = "10-2-2023"
$ws.Range("G46").Value = "0"
$rng_46_0.ClearFormats()

$rng_47_0 = $ws.Range("D47:G47")
$rng_47_0.NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000753"
$ws.Range("E47").Value = "0.38%"
$ws.Range("F47").Value = "10-2-2023"
$ws.Range("G47").Value = "0"
$rng_47_0.ClearFormats()

$rng_48_0 = $ws.Range("D48:G48")
$rng_48_0.NumberFormat = "@"
$ws.Range("D48").Value = "0.003392"
$ws.Range("E48").Value = "-2.83%"
$ws.Range("F48").Value = "10-2-2023"
$ws.Range("G48").Value = "0"
$rng_48_0.ClearFormats()

$rng_49_0 = $ws.Range("E49:G49")
$rng_49_0.NumberFormat = "@"
$ws.Range("E49").Value = "16.50%"
$ws.Range("F49").Value = "10-2-2023"
$ws.Range("G49").Value = "0"
$rng_49_0.ClearFormats()

$rng_50_0 = $ws.Range("D50:G50")
$rng_50_0.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").Value = "0.38%"
$ws.Range("F50").Value = "10-2-2023"
$ws.Range("G50").Value = "0"
$rng_50_0.ClearFormats()

$rng_51_0 = $ws.Range("D51:G51")
$rng_51_0.NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").Value = "0.38%"
$ws.Range("F51").Value = "10-2-2023"
$ws.Range("G51").Value = "0"
$rng_51_0.ClearFormats()
